$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 108.2
$ws.Range("I33").Value = 108.2
$ws.Range("K33").Value = 108.2
$ws.Range("M33").Value = 120.8

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4449
$ws.Range("I111").Value = 2863
$ws.Range("J111").Value = 10000
$ws.Range("K111").Value = 8589
$ws.Range("L111").Value = 30000
$ws.Range("M111").Value = -5522
$ws.Range("N111").Value = -36134

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2174.0588
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = ""

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3947.5
$ws.Range("J113").Value = 3900
$ws.Range("L113").Value = 3900
$ws.Range("N113").Value = -10408

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4482.4287
$ws.Range("I131").Value = 1056.75
$ws.Range("J131").Value = 9050
$ws.Range("K131").Value = 3170.25
$ws.Range("L131").Value = 27150
$ws.Range("M131").Value = 1869.75
$ws.Range("N131").Value = -37230

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1137.8788
$ws.Range("I135").Value = 1025.96
$ws.Range("K135").Value = 9233.639999999999
$ws.Range("M135").Value = -6698.639999999999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1724.1
$ws.Range("I137").Value = 1354.875
$ws.Range("K137").Value = 4064.625
$ws.Range("M137").Value = -1514.625

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3567.9534
$ws.Range("J138").Value = 2596.919
$ws.Range("L138").Value = 7790.757
$ws.Range("N138").Value = -18070.757

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5572.5537
$ws.Range("I32").Value = 3838.83
$ws.Range("K32").Value = 3838.83
$ws.Range("M32").Value = -3551.83

# ARM row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 4126.75
$ws.Range("I35").Value = 4126.75
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4126.75
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3720.75
$ws.Range("N35").Value = ""

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1262.7742
$ws.Range("I61").Value = 1163.7241
$ws.Range("K61").Value = 1163.7241
$ws.Range("M61").Value = -951.7240999999999

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1636.12
$ws.Range("I74").Value = 1025.6818
$ws.Range("K74").Value = 1025.6818
$ws.Range("M74").Value = -151.6818000000001

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1636.12
$ws.Range("I77").Value = 1025.6818
$ws.Range("K77").Value = 5128.409000000001
$ws.Range("M77").Value = -760.4090000000006

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1262.7742
$ws.Range("I136").Value = 1163.7241
$ws.Range("K136").Value = 3491.1723
$ws.Range("M136").Value = -941.1722999999997

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3061.5
$ws.Range("I20").Value = 2998.8572
$ws.Range("K20").Value = 2998.8572
$ws.Range("M20").Value = -2751.8572

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2581.0908
$ws.Range("I86").Value = 2589.2
$ws.Range("K86").Value = 2589.2
$ws.Range("M86").Value = -1466.2

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2581.0908
$ws.Range("I89").Value = 2589.2
$ws.Range("K89").Value = 12946
$ws.Range("M89").Value = -7330

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2617.5557
$ws.Range("I134").Value = 2308.7693
$ws.Range("J134").Value = 3420.4
$ws.Range("K134").Value = 6926.3079
$ws.Range("L134").Value = 10261.2
$ws.Range("M134").Value = -4391.3079
$ws.Range("N134").Value = -15331.2

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3521.2104
$ws.Range("I58").Value = 1325.4
$ws.Range("J58").Value = 5961
$ws.Range("K58").Value = 1325.4
$ws.Range("L58").Value = 5961
$ws.Range("M58").Value = -1122.4
$ws.Range("N58").Value = -6367

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12920.05
$ws.Range("I99").Value = 8654.223
$ws.Range("K99").Value = 8654.223
$ws.Range("M99").Value = -7156.223

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12920.05
$ws.Range("I126").Value = 8654.223
$ws.Range("K126").Value = 25962.669
$ws.Range("M126").Value = -23492.669

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3173.9167
$ws.Range("I134").Value = 2072.5715
$ws.Range("K134").Value = 6217.7145
$ws.Range("M134").Value = -3682.7145

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3521.2104
$ws.Range("I136").Value = 1325.4
$ws.Range("J136").Value = 5961
$ws.Range("K136").Value = 3976.2
$ws.Range("L136").Value = 17883
$ws.Range("M136").Value = -1426.2
$ws.Range("N136").Value = -22983

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2323.3333
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2323.3333
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""

# GSM row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 59195.117
$ws.Range("I3").Value = 77124.46000000001
$ws.Range("J3").Value = 924.75
$ws.Range("K3").Value = 77124.46000000001
$ws.Range("L3").Value = 924.75
$ws.Range("M3").Value = -77008.46000000001
$ws.Range("N3").Value = -1156.75

# GSM row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10025000
$ws.Range("I10").Value = 10025000
$ws.Range("K10").Value = 10025000
$ws.Range("M10").Value = -10024831

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1423800.8
$ws.Range("I11").Value = 2529500
$ws.Range("K11").Value = 2529500
$ws.Range("M11").Value = -2529361

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2310.121
$ws.Range("I102").Value = 929.6
$ws.Range("K102").Value = 929.6
$ws.Range("M102").Value = 692.4

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 114048.11
$ws.Range("I122").Value = 2928
$ws.Range("J122").Value = 202944.2
$ws.Range("K122").Value = 8784
$ws.Range("L122").Value = 608832.6000000001
$ws.Range("M122").Value = -6334
$ws.Range("N122").Value = -613732.6000000001

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5115.4165
$ws.Range("I22").Value = 783
$ws.Range("J22").Value = 6559.5557
$ws.Range("K22").Value = 783
$ws.Range("L22").Value = 6559.5557
$ws.Range("M22").Value = -488
$ws.Range("N22").Value = -7149.5557

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5115.4165
$ws.Range("I27").Value = 783
$ws.Range("J27").Value = 6559.5557
$ws.Range("K27").Value = 783
$ws.Range("L27").Value = 6559.5557
$ws.Range("M27").Value = -676
$ws.Range("N27").Value = -6773.5557

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3219.1428
$ws.Range("I68").Value = 2847.25
$ws.Range("K68").Value = 2847.25
$ws.Range("M68").Value = -2098.25

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3219.1428
$ws.Range("I71").Value = 2847.25
$ws.Range("K71").Value = 14236.25
$ws.Range("M71").Value = -10492.25

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4194.029
$ws.Range("I132").Value = 3571.1904
$ws.Range("K132").Value = 10713.5712
$ws.Range("M132").Value = -8183.5712

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 93750
$ws.Range("J133").Value = 93750
$ws.Range("L133").Value = 93750
$ws.Range("N133").Value = -98810

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7493
$ws.Range("I62").Value = 4837.2
$ws.Range("J62").Value = 8274.117
$ws.Range("K62").Value = 4837.2
$ws.Range("L62").Value = 8274.117
$ws.Range("M62").Value = -4213.2
$ws.Range("N62").Value = -9522.117

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7493
$ws.Range("I65").Value = 4837.2
$ws.Range("J65").Value = 8274.117
$ws.Range("K65").Value = 24186
$ws.Range("L65").Value = 41370.585
$ws.Range("M65").Value = -21066
$ws.Range("N65").Value = -47610.585

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 75750
$ws.Range("J133").Value = 75750
$ws.Range("L133").Value = 75750
$ws.Range("N133").Value = -85870
